# Weekly refresh of the Ciboulette price series:
# a new week's record is inserted at row 158 (pushing every existing
# record at/after that position down by one row), and the oldest
# historical record (previously the last row, 259) keeps sliding down
# to become the new final row, 260.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 158 and everything below it down by one row.
$ws.Rows.Item(158).Insert()

# Populate the freshly inserted row with the new week's data point.
$ws.Range("A158").Value = 3
$ws.Range("B158").Value = "Femacal de La Calera"
$ws.Range("C158").Value = "Coquimbo"
$ws.Range("D158").Value = 44596
$ws.Range("E158").Value = 5
$ws.Range("F158").Value = 100112039
$ws.Range("G158").Value = "Ciboulette"
$ws.Range("H158").Value = "Sin especificar"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 160
$ws.Range("K158").Value = 1500
$ws.Range("L158").Value = 1500
$ws.Range("M158").Value = 1500
$ws.Range("N158").Value = "$/docena de atados"
$ws.Range("O158").Value = "Provincia de Quillota"
$ws.Range("P158").Value = 500
$ws.Range("Q158").Value = 3
$ws.Range("R158").Value = "Hortaliza"
